$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 6
$ws.Range("S2").Value = 1.62
$ws.Range("T2").Value = 2.2
$ws.Range("AT2").Value = 2.2
$ws.Range("AU2").Value = 9.5
$ws.Range("AY2").Value = 34

# Row 8 updates
$ws.Range("G8").Value = 1.9
$ws.Range("H8").Value = 3
$ws.Range("I8").Value = 4.2
$ws.Range("J8").Value = 2.63
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 8
$ws.Range("U8").Value = 1.91
$ws.Range("V8").Value = 1.8
$ws.Range("X8").Value = 8.5
$ws.Range("AC8").Value = 7.5
$ws.Range("AJ8").Value = 41
$ws.Range("AO8").Value = 11
$ws.Range("AQ8").Value = 41
$ws.Range("AR8").Value = 67
$ws.Range("BA8").Value = 101
